$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.388.02'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.07%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.571.77'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.10%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("E5").Value = '  +0.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '291.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.05%  '

$ws.Range("E7").Value = '  +2.29%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '50.05'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3425'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07632'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.33%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.153'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.20%  '

$ws.Range("E12").Value = '  +0.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.19'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.10%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.029'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.44%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.940'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.50%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.570.76'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.17%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001132'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.72%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.08'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06747'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.36%  '

$ws.Range("E20").Value = '  +0.05%  '

$ws.Range("E21").Value = '  +2.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.197'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.56%  '

$ws.Range("E23").Value = '  -0.34%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.394.34'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.13%  '

$ws.Range("E25").Value = '  +0.13%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.679'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -10.75%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.19'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.47%  '

$ws.Range("E28").Value = '  +1.24%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.029'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.86%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.24'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.60%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.746.06'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.142'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.91%  '

$ws.Range("E33").Value = '  +0.75%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9841'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.36%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.923'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.32%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08540'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.98%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02546'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.13%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2316'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.55%  '

$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.341'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.22%  '

$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06562'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.23%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.413'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.22%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6396'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.42'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.55%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.04'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.15%  '

$ws.Range("E46").Value = '  +0.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5982'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.28%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.303'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.68%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.086'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.55%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '125.12'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.45%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07324'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.34%  '
